$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose only change is column B: 98926 -> 98930
$simpleRows = @(2,3,4,5,6,7,8,9,10,13,14)
foreach ($r in $simpleRows) {
    $ws.Range("B$r").Value = 98930
}

# Rows 11 and 12 swap their record contents (row 11 <-> row 12),
# with the Taxonsorteringsordning (column B) values updated to the new ones.

# New row 11 (was the NT/Garnlav record previously on row 12, with B updated 79239 -> 79243)
$ws.Range("A11").Value = 130965935
$ws.Range("B11").Value = 79243
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 496969
$ws.Range("R11").Value = 6713674
$ws.Range("AC11").Value = "Måttlig förekomst . inventering åt vasa vind"
$ws.Range("AX11").Value = "Pia Edfors, Enviro Planning"

# New row 12 (was the LC/Fläcknycklar record previously on row 11, with B updated 98926 -> 98930)
$ws.Range("A12").Value = 130965861
$ws.Range("B12").Value = 98930
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 219790
$ws.Range("F12").Value = "Fläcknycklar"
$ws.Range("G12").Value = "Dactylorhiza maculata"
$ws.Range("H12").Value = "(L.) Soó"
$ws.Range("Q12").Value = 497138
$ws.Range("R12").Value = 6713448
$ws.Range("AC12").Value = "Betydelsefulla förekomster . inventering åt vasa vind"
$ws.Range("AX12").Value = "Anders Esplund, Pia Edfors, Enviro Planning"
